$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "products" template now also needs to carry the photo attachment
# file names, so add the two extra lookup columns to the header row.
$ws.Range("K1").Value2 = "file_name"
$ws.Range("L1").Value2 = "tile_file_name"

# Re-fit the columns so the new/changed header text is fully visible
# (mirrors doing Home > Format > AutoFit Column Width on these columns).
$ws.Columns.Item(1).EntireColumn.AutoFit()  | Out-Null   # A - product_code
$ws.Columns.Item(4).EntireColumn.AutoFit()  | Out-Null   # D - brand
$ws.Columns.Item(9).EntireColumn.AutoFit()  | Out-Null   # I - price
$ws.Columns.Item(10).EntireColumn.AutoFit() | Out-Null   # J - range
$ws.Columns.Item(11).EntireColumn.AutoFit() | Out-Null   # K - file_name
$ws.Columns.Item(12).EntireColumn.AutoFit() | Out-Null   # L - tile_file_name

# Leave the selection where the user ended up after the edit.
$ws.Range("C8").Select() | Out-Null
